# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de), the handback process has now
# completed for both rows. This fills in the "Latest Target File" (F) and
# "Latest Handback File" (G) columns (mirroring the source name in A and the
# handoff xlf in D), stamps the "Latest Handback DateTime" (H) and updates
# the "Status" (C) to reflect the sync.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkForAddress($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            return $h
        }
    }
    return $null
}

function Add-MirrorHyperlink($ws, $srcAddr, $destCell, $destText) {
    $src = Get-HyperlinkForAddress $ws $srcAddr
    if ($src -ne $null) {
        $ws.Hyperlinks.Add($ws.Range($destCell), $src.Address(), "", "", $destText) | Out-Null
    } else {
        $ws.Range($destCell).Value = $destText
    }
}

function Update-HandbackSheet($ws) {
    # Row 2 -- 5e72deef-91dd-4901-bf8e-bcc72ab6d1f4
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    Add-MirrorHyperlink $ws "`$A`$2" "F2" $ws.Range("A2").Value()
    Add-MirrorHyperlink $ws "`$D`$2" "G2" $ws.Range("D2").Value()
    $ws.Range("H2").Value = "2016-03-20 02:12:38"

    # Row 3 -- 9a7d9abc-9216-491f-a1df-57ef41339a1b
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
    Add-MirrorHyperlink $ws "`$A`$3" "F3" $ws.Range("A3").Value()
    Add-MirrorHyperlink $ws "`$D`$3" "G3" $ws.Range("D3").Value()
    $ws.Range("H3").Value = "2016-03-20 02:12:43"
}

Update-HandbackSheet $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wb.Worksheets.Item("de-de")
